# Updated cryptos list on Mon Oct 16 22:00:00 UTC 2023 with GitHub Actions
# This script updates the Price (D) and Volume(1h) (E) columns for the crypto
# ranking sheet, and reorders a handful of rows whose relative ranking changed
# (MXToken/ImmutableX, WEMIXToken/RenderToken, and the new mCoin entry pushing
# BabyDogeCoin/Cronos down and dropping BitcoinSV off the bottom of the list).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.430.77"
$ws.Range("E2").Value = "  +4.53%  "

$ws.Range("D3").Value = "1.590.76"
$ws.Range("E3").Value = "  +1.84%  "

$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").Value = "'214.65"
$ws.Range("E5").Value = "  +2.16%  "

$ws.Range("D6").Value = "'0.498"
$ws.Range("E6").Value = "  +1.48%  "

$ws.Range("E7").Value = "  -0.05%  "

$ws.Range("D8").Value = "'23.97"
$ws.Range("E8").Value = "  +8.88%  "

$ws.Range("E9").Value = "  +1.25%  "

$ws.Range("E10").Value = "  +0.78%  "

$ws.Range("E11").Value = "  +2.40%  "

$ws.Range("D12").Value = "1.817.44"
$ws.Range("E12").Value = "  +1.78%  "

$ws.Range("D13").Value = "1.579.81"
$ws.Range("E13").Value = "  +1.17%  "

$ws.Range("D14").Value = "'3.80"
$ws.Range("E14").Value = "  +0.83%  "

$ws.Range("E15").Value = "  +3.19%  "

$ws.Range("D16").Value = "28.443.46"
$ws.Range("E16").Value = "  +4.72%  "

$ws.Range("D17").Value = "'63.13"
$ws.Range("E17").Value = "  +2.05%  "

$ws.Range("D18").Value = "'232.49"
$ws.Range("E18").Value = "  +7.30%  "

$ws.Range("E19").Value = "  +1.20%  "

$ws.Range("E20").Value = "  +0.51%  "

$ws.Range("E21").Value = "  -0.12%  "

$ws.Range("D22").Value = "'4.13"
$ws.Range("E22").Value = "  -0.13%  "

$ws.Range("D23").Value = "'9.44"
$ws.Range("E23").Value = "  +2.86%  "

$ws.Range("E24").Value = "  +1.17%  "

$ws.Range("D25").Value = "'152.16"
$ws.Range("E25").Value = "  -0.40%  "

$ws.Range("D26").Value = "'15.30"
$ws.Range("E26").Value = "  +2.02%  "

$ws.Range("E27").Value = "  -0.09%  "

$ws.Range("D28").Value = "'0.108"
$ws.Range("E28").Value = "  +1.44%  "

$ws.Range("E29").Value = "  -0.11%  "

$ws.Range("E30").Value = "  +0.92%  "

$ws.Range("E31").Value = "  +1.01%  "

$ws.Range("E32").Value = "  +0.54%  "

$ws.Range("E33").Value = "  +0.45%  "

$ws.Range("D34").Value = "1.420.48"
$ws.Range("E34").Value = "  -1.18%  "

$ws.Range("E35").Value = "  -0.72%  "

$ws.Range("E36").Value = "  -4.63%  "

$ws.Range("E37").Value = "  -0.11%  "

$ws.Range("E38").Value = "  +0.69%  "

$ws.Range("B39").Value = "ImmutableX"
$ws.Range("C39").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D39").Value = "'0.546"
$ws.Range("E39").Value = "  +2.53%  "

$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").Value = "'2.53"
$ws.Range("E40").Value = "  +5.48%  "

$ws.Range("D41").Value = "'0.821"
$ws.Range("E41").Value = "  +1.93%  "

$ws.Range("E42").Value = "  -2.58%  "

$ws.Range("E43").Value = "  -0.13%  "

$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").Value = "'1.83"
$ws.Range("E44").Value = "  +6.18%  "

$ws.Range("B45").Value = "WEMIXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D45").Value = "'0.978"
$ws.Range("E45").Value = "  -2.11%  "

$ws.Range("D46").Value = "'64.73"
$ws.Range("E46").Value = "  +0.72%  "

$ws.Range("D47").Value = "1.729.89"
$ws.Range("E47").Value = "  +1.62%  "

$ws.Range("D48").Value = "'87.68"
$ws.Range("E48").Value = "  +2.18%  "

$ws.Range("B49").Value = "mCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin"
$ws.Range("D49").Value = "'2.14"
$ws.Range("E49").Value = "  +0.50%  "

$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "0.0₆0108"
$ws.Range("E50").Value = "  +12.03%  "

$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "'0.0523"
$ws.Range("E51").Value = "  -0.41%  "
